$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.301.73'
$ws.Range('E2').Value = '  +2.12%  '
$ws.Range('D3').Value = '1.658.00'
$ws.Range('E3').Value = '  +1.17%  '
$ws.Range('E4').Value = '  -0.55%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '219.45'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.92%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.506'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.60%  '
$ws.Range('E7').Value = '  -0.51%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.257'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.46%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0626'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.47%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.98'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.81%  '
$ws.Range('E11').Value = '  +0.44%  '
$ws.Range('D12').Value = '1.890.18'
$ws.Range('E12').Value = '  +1.16%  '
$ws.Range('D13').Value = '1.670.85'
$ws.Range('E13').Value = '  +2.19%  '
$ws.Range('E14').Value = '  +1.21%  '
$ws.Range('E15').Value = '  +1.09%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.21'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.34%  '
$ws.Range('D17').Value = '27.290.21'
$ws.Range('E17').Value = '  +2.10%  '
$ws.Range('D18').Value = '0.0₃0735'
$ws.Range('E18').Value = '  +0.89%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '222.38'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.96%  '
$ws.Range('E20').Value = '  -0.49%  '
$ws.Range('E21').Value = '  +2.03%  '
$ws.Range('E22').Value = '  +8.21%  '
$ws.Range('E23').Value = '  +4.49%  '
$ws.Range('E24').Value = '  +0.30%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.18'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.13%  '
$ws.Range('E26').Value = '  -0.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.42'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.38%  '
$ws.Range('E28').Value = '  +2.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '16.04'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.82%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0514'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.73%  '
$ws.Range('E31').Value = '  +0.84%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.40'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.88%  '
$ws.Range('E33').Value = '  +0.66%  '
$ws.Range('E34').Value = '  +2.47%  '
$ws.Range('D35').Value = '1.262.11'
$ws.Range('E35').Value = '  -1.15%  '
$ws.Range('E36').Value = '  +1.08%  '
$ws.Range('E37').Value = '  +1.90%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.538'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.836'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.93%  '
$ws.Range('E40').Value = '  -0.47%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.818'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.66%  '
$ws.Range('E42').Value = '  +2.10%  '
$ws.Range('D43').Value = '1.800.47'
$ws.Range('E43').Value = '  +1.28%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.12'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.85%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '61.78'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.62%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '91.98'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.01%  '
$ws.Range('E47').Value = '  +1.56%  '
$ws.Range('E48').Value = '  -0.91%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0980'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.84%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.65'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.37%  '
$ws.Range('E51').Value = '  +0.16%  '
